# Apply the "Web page" slide (solution slides) box-added edit:
#  - Add a new accent1-styled rectangle ("Rectangle 3") behind the existing
#    content on slide 16 ("Web page").
#  - Reposition / resize the four picture icons that sit on top of it.
#
# Note: Length properties (Left/Top/Width/Height) are expressed in points in
# the PowerPoint object model, while the underlying OOXML stores EMU
# (1 pt = 12700 EMU). To land exactly on the target EMU values despite the
# point<->EMU round trip, a tiny (+0.5 EMU) epsilon is added before dividing.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(16)

function EmuToPt($emu) {
    return ($emu + 0.5) / 12700
}

# --- 1. Add the new rectangle behind everything else -----------------------
# There is no direct "insert default-styled autoshape" COM call available, so
# borrow the style (p:style + plain spPr, no explicit fill/line overrides)
# from an existing plain accent1 shape elsewhere in the deck, turn it into a
# plain rectangle, then move/resize/rename it.
$styleDonor = $p.Slides.Item(8).Shapes.Item(9)
$styleDonor.Copy()
$pasted = $slide.Shapes.Paste()
$rect = $pasted.Item(1)

$rect.AutoShapeType = 1        # msoShapeRectangle
$rect.Name = "Rectangle 3"

$rect.Left = EmuToPt 4751033
$rect.Top = EmuToPt 2898395
$rect.Width = EmuToPt 2689934
$rect.Height = EmuToPt 3719744

# Send the new rectangle to the back of the z-order (first shape in spTree).
$rect.ZOrder(1)

# --- 2. Reposition / resize the four picture icons --------------------------
$pic1 = $slide.Shapes.Item("Graphique 4")
$pic1.Left = EmuToPt 4843568
$pic1.Top = EmuToPt 3067031
$pic1.Width = EmuToPt 641781
$pic1.Height = EmuToPt 641781

$pic2 = $slide.Shapes.Item("Graphique 5")
$pic2.Left = EmuToPt 5485349
$pic2.Top = EmuToPt 3843867

$pic3 = $slide.Shapes.Item("Graphique 6")
$pic3.Left = EmuToPt 5963466
$pic3.Top = EmuToPt 4620703

$pic4 = $slide.Shapes.Item("Graphique 7")
$pic4.Left = EmuToPt 6499841
$pic4.Top = EmuToPt 3843867
